# Applies the targeted corrections to the "company_list" sheet (rows 2-9,
# columns D:AJ) as described in the commit "error solve ifrs list".
# Numeric metrics are rewritten with corrected figures and a handful of
# now-obsolete columns (J, O, and the stray V4) are cleared entirely so the
# cell no longer exists in the sheet (matching the upstream removal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3014
$ws.Range("E2").Value = 121
$ws.Range("F2").Value = 121
$ws.Range("G2").Value = 174
$ws.Range("H2").Value = 140
$ws.Range("I2").Value = 140
$ws.Range("K2").Value = 3227
$ws.Range("L2").Value = 1331
$ws.Range("M2").Value = 1895
$ws.Range("N2").Value = 1895
$ws.Range("P2").Value = 77
$ws.Range("Q2").Value = -186
$ws.Range("R2").Value = -5
$ws.Range("S2").Value = 192
$ws.Range("T2").Value = 154
$ws.Range("U2").Value = -340
$ws.Range("V2").Value = 684
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 4.64
$ws.Range("Y2").Value = 7.62
$ws.Range("Z2").Value = 4.66
$ws.Range("AA2").Value = 70.25
$ws.Range("AB2").Value = 2361.22
$ws.Range("AC2").Value = 908
$ws.Range("AD2").Value = 16.52
$ws.Range("AE2").Value = 12307
$ws.Range("AF2").Value = 1.22
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 16.52
$ws.Range("AJ2").Value = 15400000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 3700
$ws.Range("E3").Value = 188
$ws.Range("F3").Value = 188
$ws.Range("G3").Value = 168
$ws.Range("H3").Value = 124
$ws.Range("I3").Value = 124
$ws.Range("K3").Value = 3334
$ws.Range("L3").Value = 1338
$ws.Range("M3").Value = 1996
$ws.Range("N3").Value = 1996
$ws.Range("P3").Value = 77
$ws.Range("Q3").Value = 157
$ws.Range("R3").Value = -46
$ws.Range("S3").Value = -112
$ws.Range("T3").Value = 65
$ws.Range("U3").Value = 91
$ws.Range("V3").Value = 596
$ws.Range("W3").Value = 5.07
$ws.Range("X3").Value = 3.36
$ws.Range("Y3").Value = 6.39
$ws.Range("Z3").Value = 3.79
$ws.Range("AA3").Value = 67.02
$ws.Range("AB3").Value = 2492.53
$ws.Range("AC3").Value = 808
$ws.Range("AD3").Value = 21.97
$ws.Range("AE3").Value = 12964
$ws.Range("AF3").Value = 1.37
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 1.13
$ws.Range("AI3").Value = 24.76
$ws.Range("AJ3").Value = 15400000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 4390
$ws.Range("E4").Value = 456
$ws.Range("F4").Value = 456
$ws.Range("G4").Value = 411
$ws.Range("H4").Value = 302
$ws.Range("I4").Value = 302
$ws.Range("K4").Value = 3056
$ws.Range("L4").Value = 789
$ws.Range("M4").Value = 2267
$ws.Range("N4").Value = 2267
$ws.Range("P4").Value = 77
$ws.Range("Q4").Value = 558
$ws.Range("R4").Value = 167
$ws.Range("S4").Value = -626
$ws.Range("T4").Value = 42
$ws.Range("U4").Value = 516
$ws.Range("W4").Value = 10.39
$ws.Range("X4").Value = 6.89
$ws.Range("Y4").Value = 14.18
$ws.Range("Z4").Value = 9.46
$ws.Range("AA4").Value = 34.79
$ws.Range("AB4").Value = 2844.24
$ws.Range("AC4").Value = 1963
$ws.Range("AD4").Value = 8.89
$ws.Range("AE4").Value = 14721
$ws.Range("AF4").Value = 1.19
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 2.01
$ws.Range("AI4").Value = 17.83
$ws.Range("AJ4").Value = 15400000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()

# Row 5
$ws.Range("D5").Value = 5605
$ws.Range("E5").Value = 981
$ws.Range("F5").Value = 981
$ws.Range("G5").Value = 995
$ws.Range("H5").Value = 749
$ws.Range("I5").Value = 749
$ws.Range("K5").Value = 4001
$ws.Range("L5").Value = 1041
$ws.Range("M5").Value = 2960
$ws.Range("N5").Value = 2960
$ws.Range("P5").Value = 77
$ws.Range("Q5").Value = 666
$ws.Range("R5").Value = -388
$ws.Range("S5").Value = -54
$ws.Range("T5").Value = 108
$ws.Range("U5").Value = 558
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 17.5
$ws.Range("X5").Value = 13.36
$ws.Range("Y5").Value = 28.66
$ws.Range("Z5").Value = 21.22
$ws.Range("AA5").Value = 35.19
$ws.Range("AB5").Value = 3745.75
$ws.Range("AC5").Value = 4863
$ws.Range("AD5").Value = 8.63
$ws.Range("AE5").Value = 19218
$ws.Range("AF5").Value = 2.18
$ws.Range("AG5").Value = 650
$ws.Range("AH5").Value = 1.55
$ws.Range("AI5").Value = 13.37
$ws.Range("AJ5").Value = 15400000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 6687
$ws.Range("E6").Value = 915
$ws.Range("F6").Value = 915
$ws.Range("G6").Value = 1410
$ws.Range("H6").Value = 1093
$ws.Range("I6").Value = 1092
$ws.Range("K6").Value = 4958
$ws.Range("L6").Value = 1010
$ws.Range("M6").Value = 3948
$ws.Range("N6").Value = 3941
$ws.Range("P6").Value = 77
$ws.Range("Q6").Value = 273
$ws.Range("R6").Value = -330
$ws.Range("S6").Value = -120
$ws.Range("T6").Value = 477
$ws.Range("U6").Value = -204
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 13.69
$ws.Range("X6").Value = 16.34
$ws.Range("Y6").Value = 31.64
$ws.Range("Z6").Value = 24.4
$ws.Range("AA6").Value = 25.59
$ws.Range("AB6").Value = 5020.1
$ws.Range("AC6").Value = 7089
$ws.Range("AD6").Value = 5.65
$ws.Range("AE6").Value = 25593
$ws.Range("AF6").Value = 1.56
$ws.Range("AG6").Value = 650
$ws.Range("AH6").Value = 1.62
$ws.Range("AI6").Value = 9.17
$ws.Range("AJ6").Value = 15400000

# Row 7
$ws.Range("D7").Value = 9098
$ws.Range("E7").Value = 1366
$ws.Range("G7").Value = 1388
$ws.Range("H7").Value = 1046
$ws.Range("I7").Value = 1046
$ws.Range("K7").Value = 6390
$ws.Range("L7").Value = 1527
$ws.Range("M7").Value = 4864
$ws.Range("N7").Value = 4863
$ws.Range("P7").Value = 79
$ws.Range("Q7").Value = 967
$ws.Range("R7").Value = -289
$ws.Range("S7").Value = -196
$ws.Range("T7").Value = 287
$ws.Range("U7").Value = 793
$ws.Range("W7").Value = 15.02
$ws.Range("X7").Value = 11.49
$ws.Range("Y7").Value = 23.76
$ws.Range("Z7").Value = 18.43
$ws.Range("AA7").Value = 31.39
$ws.Range("AC7").Value = 6790
$ws.Range("AD7").Value = 14.95
$ws.Range("AE7").Value = 31704
$ws.Range("AF7").Value = 3.2
$ws.Range("AG7").Value = 655
$ws.Range("AH7").Value = 0.65
$ws.Range("AI7").Value = 9.65

# Row 8
$ws.Range("D8").Value = 10722
$ws.Range("E8").Value = 1672
$ws.Range("G8").Value = 1698
$ws.Range("H8").Value = 1294
$ws.Range("I8").Value = 1294
$ws.Range("K8").Value = 7697
$ws.Range("L8").Value = 1644
$ws.Range("M8").Value = 6052
$ws.Range("N8").Value = 6056
$ws.Range("P8").Value = 79
$ws.Range("Q8").Value = 1202
$ws.Range("R8").Value = -326
$ws.Range("S8").Value = -96
$ws.Range("T8").Value = 301
$ws.Range("U8").Value = 989
$ws.Range("W8").Value = 15.59
$ws.Range("X8").Value = 12.07
$ws.Range("Y8").Value = 23.69
$ws.Range("Z8").Value = 18.37
$ws.Range("AA8").Value = 27.17
$ws.Range("AC8").Value = 8399
$ws.Range("AD8").Value = 12.08
$ws.Range("AE8").Value = 39511
$ws.Range("AF8").Value = 2.57
$ws.Range("AG8").Value = 665
$ws.Range("AH8").Value = 0.66
$ws.Range("AI8").Value = 7.92

# Row 9
$ws.Range("D9").Value = 12271
$ws.Range("E9").Value = 1953
$ws.Range("G9").Value = 1991
$ws.Range("H9").Value = 1525
$ws.Range("I9").Value = 1523
$ws.Range("K9").Value = 9162
$ws.Range("L9").Value = 1676
$ws.Range("M9").Value = 7486
$ws.Range("N9").Value = 7481
$ws.Range("P9").Value = 79
$ws.Range("Q9").Value = 1519
$ws.Range("R9").Value = -478
$ws.Range("S9").Value = -96
$ws.Range("T9").Value = 409
$ws.Range("U9").Value = 1245
$ws.Range("W9").Value = 15.92
$ws.Range("X9").Value = 12.43
$ws.Range("Y9").Value = 22.5
$ws.Range("Z9").Value = 18.09
$ws.Range("AA9").Value = 22.39
$ws.Range("AC9").Value = 9890
$ws.Range("AD9").Value = 10.26
$ws.Range("AE9").Value = 48813
$ws.Range("AF9").Value = 2.08
$ws.Range("AG9").Value = 661
$ws.Range("AH9").Value = 0.65
$ws.Range("AI9").Value = 6.69
